# Update homeroom/admin teacher assignments for "Chào cờ" / "SHL" / "HĐNGLL"
# rows in both "Học kì 1" (rows 20-22) and "Học kì 2" (rows 51-53) blocks of
# each class sheet. Previously all of these used a single placeholder value
# ("homeroomteacher"); now "Chào cờ" is taught by "admin" while "SHL" and
# "HĐNGLL" are taught by a per-class homeroom teacher.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("12A1")
$ws1.Range("C20").Value = "admin"
$ws1.Range("C21").Value = "homeroomteacher12"
$ws1.Range("C22").Value = "homeroomteacher12"
$ws1.Range("C51").Value = "admin"
$ws1.Range("C52").Value = "homeroomteacher12"
$ws1.Range("C53").Value = "homeroomteacher12"

$ws2 = $wb.Worksheets.Item("12A2")
$ws2.Range("C20").Value = "admin"
$ws2.Range("C21").Value = "homeroomteacher23"
$ws2.Range("C22").Value = "homeroomteacher23"
$ws2.Range("C51").Value = "admin"
$ws2.Range("C52").Value = "homeroomteacher23"
$ws2.Range("C53").Value = "homeroomteacher23"

$ws3 = $wb.Worksheets.Item("12A3")
$ws3.Range("C20").Value = "admin"
$ws3.Range("C21").Value = "homeroomteacher34"
$ws3.Range("C22").Value = "homeroomteacher34"
$ws3.Range("C51").Value = "admin"
$ws3.Range("C52").Value = "homeroomteacher34"
$ws3.Range("C53").Value = "homeroomteacher34"

# Reproduce the final selection/scroll state recorded in the workbook: the
# user last worked on 12A2 and 12A3 around the C51:C53 block, then ended up
# back on 12A1 (which stays the active/displayed tab) with H17 selected.
$ws2.Activate()
$ws2.Range("C51:C53").Select()

$ws3.Activate()
$ws3.Range("C51:C53").Select()

$ws1.Activate()
$ws1.Range("H17").Select()
